$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 3
    "I2" = 2.2
    "J2" = 3.4
    "L2" = 2.75
    "N2" = 17
    "U2" = 1.47
    "X2" = 19
    "Z2" = 34
    "AE2" = 12
    "AH2" = 11
    "AJ2" = 9
    "AK2" = 21
    "AL2" = 15
    "AN2" = 5.5
    "AP2" = 21
    "AQ2" = 51
    "AX2" = 4.5
    "AY2" = 11
    "BA2" = 34
    "U3" = 1.58
    "BD3" = 176
    "J4" = 1.5
    "K4" = 3.1
    "M4" = 1.02
    "N4" = 21
    "O4" = 1.11
    "P4" = 6.5
    "Q4" = 1.4
    "R4" = 2.88
    "S4" = 1.22
    "T4" = 4
    "U4" = 2.2
    "V4" = 1.58
    "W4" = 9.5
    "X4" = 6.5
    "AD4" = 15
    "AG4" = 501
    "AM4" = 81
    "AR4" = 34
    "AS4" = 126
    "AT4" = 4
    "BA4" = 401
    "BC4" = 451
    "G5" = 4.33
    "H5" = 3.6
    "I5" = 1.85
    "J5" = 4.5
    "L5" = 2.5
    "Q5" = 1.89
    "R5" = 1.84
    "U5" = 1.77
    "V5" = 1.92
    "X5" = 21
    "Z5" = 41
    "AC5" = 10
    "AD5" = 6.5
    "AH5" = 7.5
    "AY5" = 10
    "BA5" = 34
    "G6" = 4.33
    "I6" = 1.73
    "M6" = 1.03
    "N6" = 17
    "Q6" = 1.5
    "R6" = 2.63
    "U6" = 1.41
    "V6" = 2.62
    "AP6" = 21
    "AX6" = 4.33
    "BD6" = 151
    "G7" = 1.55
    "H7" = 4
    "I7" = 5.5
    "J7" = 2.1
    "L7" = 5.5
    "M7" = 1.02
    "O7" = 1.17
    "X7" = 8
    "Y7" = 8
    "AD7" = 8
    "AN7" = 3.6
    "AS7" = 101
    "BA7" = 101
    "H8" = 3.4
    "I8" = 3
    "L8" = 3.6
    "N8" = 12
    "W8" = 9
    "AB8" = 23
    "AC8" = 12
    "AG8" = 151
    "AJ8" = 11
    "AS8" = 126
    "G10" = 2.82
    "I10" = 2.55
    "G11" = 1.91
    "G12" = 1.79
    "G13" = 2.32
    "I13" = 2.75
    "G15" = 1.73
    "H15" = 4
    "I15" = 4.33
    "J15" = 2.25
    "R15" = 2.5
    "U15" = 1.53
    "V15" = 2.38
    "X15" = 10
    "AA15" = 12
    "AC15" = 17
    "AI15" = 26
    "AJ15" = 15
    "AO15" = 8.5
    "AU15" = 7.5
    "BA15" = 67
    "G17" = 1.55
    "H17" = 4.33
    "I17" = 5.25
    "W17" = 7.5
    "AG17" = 251
    "AJ17" = 17
    "AY17" = 29
    "AZ17" = 34
    "H18" = 2.9
    "I18" = 3.5
    "K18" = 1.91
    "L18" = 4.33
    "O18" = 1.5
    "P18" = 2.5
    "Q18" = 2.6
    "R18" = 1.48
    "S18" = 1.57
    "T18" = 2.25
    "U18" = 2.2
    "V18" = 1.62
    "AC18" = 6.5
    "AE18" = 19
    "AS18" = 301
    "AT18" = 2.25
    "AU18" = 9.5
    "AV18" = 81
    "BA18" = 81
    "BB18" = 126
    "BC18" = 351
    "I19" = 1.44
    "L19" = 1.91
    "Q19" = 1.53
    "R19" = 2.32
    "N20" = 9.5
    "O20" = 1.33
    "P20" = 3.25
    "R20" = 1.72
    "AB20" = 34
    "G21" = 1.96
    "I21" = 3.8
    "J21" = 2.63
    "L21" = 4
    "Q21" = 1.9
    "R21" = 1.95
    "X21" = 10
    "Y21" = 9
    "AB21" = 26
    "AI21" = 19
    "AO21" = 11
    "AX21" = 5.5
    "AZ21" = 26
    "BC21" = 151
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
